$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.338.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5318"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2655"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07867"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.534"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.665.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.900.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5620"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.361.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.728"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "200.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.39%  "

$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.066"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.012"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1216"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.266"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05902"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.287"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.524"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.332"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9667"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.831"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.433"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.971"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.077.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8616"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.810.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4420"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05152"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
